$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing measurement values (wooden panels measurements)
$ws.Range("B1").Value = 267
$ws.Range("B3").Value = 630
$ws.Range("B6").Value = 183
$ws.Range("B7").Value = 53

# Add new rows 19-21
$ws.Range("A19").Value = "w_infosign"
$ws.Range("B19").Value = 100
$ws.Range("C19").Value = "mm"

$ws.Range("A20").Value = "l_infosign"
$ws.Range("B20").Value = 180
$ws.Range("C20").Value = "mm"

$ws.Range("A21").Value = "w_hidden_panel"
$ws.Range("B21").Value = 75
$ws.Range("C21").Value = "mm"
$ws.Range("D21").Value = "Part that is hidden by another piece"

# Restore selection to B4 as in the final saved state
$ws.Range("B4").Select()
